$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A1").Value = "hello world"
Write-Host "Sheet1 name: $($ws.Name)"
Write-Host "Value A1: $($ws.Range('A1').Value)"
